$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 (old "A/AO1 Pin2 / VEE") is removed; rows 3-18 shift up by one ----
$ws.Rows.Item(2).Delete()

# ---- A new row is inserted at 18 (pushes old rows 18-23, now 18-23, down to 19-24) ----
$ws.Rows.Item(18).Insert(-4121)

# ---- Give the new row-18 cells the same plain bordered style as its neighbours ----
$ws.Range("B19").Copy()
$ws.Range("B18:C18").PasteSpecial(-4122)

# ---- New "Moteurs sur plaque" block: columns E-G, rows 5-11 ----
$ws.Range("B1").Copy()
$ws.Range("E5:G5").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("F6:G11").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("E6:E11").PasteSpecial(-4122)
$ws.Range("E6:E11").Interior.Color = 49407   # RGB(255,192,0) = FFC000 (orange)

# ---- Cell values (text overwrites + the brand new cells) ----
$ws.Range("B2").Value = "A/AO1 (Pin 4)"
$ws.Range("C2").Value = "Demux: C"
$ws.Range("B17").Value = "C/AO0 (AO0)"
$ws.Range("C17").Value = "VDD (cst = 5)"
$ws.Range("B18").Value = "C/AO1 (AO1)"
$ws.Range("C18").Value = "VEE (cst = -10)"

$ws.Range("E5").Value = "Moteurs sur plaque:"
$ws.Range("F5").Value = "Driver"
$ws.Range("G5").Value = "Encodeur"
$ws.Range("E6").Value = "F"
$ws.Range("F6").Value = "M1"
$ws.Range("G6").Value = "Enco1"
$ws.Range("E7").Value = "R"
$ws.Range("F7").Value = "M2"
$ws.Range("G7").Value = "Enco2"
$ws.Range("E8").Value = "B"
$ws.Range("F8").Value = "M3"
$ws.Range("G8").Value = "Enco3"
$ws.Range("E9").Value = "L"
$ws.Range("F9").Value = "M4"
$ws.Range("G9").Value = "Enco6"
$ws.Range("E10").Value = "D"
$ws.Range("F10").Value = "M5"
$ws.Range("G10").Value = "Enco4"
$ws.Range("E11").Value = "U"
$ws.Range("F11").Value = "M6"
$ws.Range("G11").Value = "Enco5"

# ---- New column E width ----
$ws.Columns.Item(5).ColumnWidth = 17.92

# ---- Selection ----
$ws.Range("D3").Select()
